$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.03681468963623
$ws.Range("B1").Value = 1.927991390228271
$ws.Range("C1").Value = 8.018400192260742
$ws.Range("D1").Value = 1.815383791923523
$ws.Range("E1").Value = 0.7094436287879944
